$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 202498.9
$ws.Range("I46").Value = 9500
$ws.Range("J46").Value = 250748.62
$ws.Range("K46").Value = 28500
$ws.Range("L46").Value = 752245.86
$ws.Range("M46").Value = -28381
$ws.Range("N46").Value = -752483.86
$ws.Range("H53").Value = 2455.3845
$ws.Range("I53").Value = 427.8
$ws.Range("J53").Value = 3722.625
$ws.Range("K53").Value = 427.8
$ws.Range("L53").Value = 3722.625
$ws.Range("M53").Value = 209.2
$ws.Range("N53").Value = -4996.625
$ws.Range("H60").Value = 202498.9
$ws.Range("I60").Value = 9500
$ws.Range("J60").Value = 250748.62
$ws.Range("K60").Value = 28500
$ws.Range("L60").Value = 752245.86
$ws.Range("M60").Value = -28016
$ws.Range("N60").Value = -753213.86
$ws.Range("H62").Value = 7749.5
$ws.Range("J62").Value = 11999
$ws.Range("L62").Value = 11999
$ws.Range("N62").Value = -13247
$ws.Range("H65").Value = 7749.5
$ws.Range("J65").Value = 11999
$ws.Range("L65").Value = 59995
$ws.Range("N65").Value = -66235
$ws.Range("H137").Value = 5994.4116
$ws.Range("I137").Value = 5464.643
$ws.Range("K137").Value = 16393.929
$ws.Range("M137").Value = -13843.929
$ws.Range("H138").Value = 107379.414
$ws.Range("J138").Value = 119373.07
$ws.Range("L138").Value = 358119.21
$ws.Range("N138").Value = -368399.21

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 3599.4
$ws.Range("I4").Value = 2332.3333
$ws.Range("K4").Value = 2332.3333
$ws.Range("M4").Value = -2216.3333
$ws.Range("H32").Value = 12782.159
$ws.Range("I32").Value = 12023.375
$ws.Range("K32").Value = 12023.375
$ws.Range("M32").Value = -11736.375
$ws.Range("H45").Value = 2788.7778
$ws.Range("J45").Value = 2880.6
$ws.Range("L45").Value = 2880.6
$ws.Range("N45").Value = -3634.6
$ws.Range("H110").Value = 1152.5
$ws.Range("I110").Value = 934.8461
$ws.Range("K110").Value = 934.8461
$ws.Range("M110").Value = 1110.1539
$ws.Range("H132").Value = 1662.1531
$ws.Range("I132").Value = 1475.9512
$ws.Range("J132").Value = 2616.4375
$ws.Range("K132").Value = 4427.8536
$ws.Range("L132").Value = 7849.3125
$ws.Range("M132").Value = -1897.8536
$ws.Range("N132").Value = -12909.3125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3439.037
$ws.Range("I20").Value = 3374.238
$ws.Range("J20").Value = 3665.8333
$ws.Range("K20").Value = 3374.238
$ws.Range("L20").Value = 3665.8333
$ws.Range("M20").Value = -3127.238
$ws.Range("N20").Value = -4159.8333
$ws.Range("H86").Value = 3036035.8
$ws.Range("I86").Value = 4448177.5
$ws.Range("K86").Value = 4448177.5
$ws.Range("M86").Value = -4447054.5
$ws.Range("H89").Value = 3036035.8
$ws.Range("I89").Value = 4448177.5
$ws.Range("K89").Value = 22240887.5
$ws.Range("M89").Value = -22235271.5
$ws.Range("H94").Value = 2899.1724
$ws.Range("J94").Value = 7811.857
$ws.Range("L94").Value = 7811.857
$ws.Range("N94").Value = -8713.857
$ws.Range("H134").Value = 2755.7693
$ws.Range("I134").Value = 2560.2952
$ws.Range("K134").Value = 7680.8856
$ws.Range("M134").Value = -5145.8856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 768.2381
$ws.Range("I7").Value = 143.06667
$ws.Range("J7").Value = 2331.1667
$ws.Range("K7").Value = 143.06667
$ws.Range("L7").Value = 2331.1667
$ws.Range("M7").Value = -30.06666999999999
$ws.Range("N7").Value = -2557.1667
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H60").Value = 55871.285
$ws.Range("J60").Value = 62220
$ws.Range("L60").Value = 62220
$ws.Range("N60").Value = -63242
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 43573.75
$ws.Range("I68").Value = 43000
$ws.Range("J68").Value = 43765
$ws.Range("K68").Value = 43000
$ws.Range("L68").Value = 43765
$ws.Range("M68").Value = -42251
$ws.Range("N68").Value = -45263
$ws.Range("H71").Value = 43573.75
$ws.Range("I71").Value = 43000
$ws.Range("J71").Value = 43765
$ws.Range("K71").Value = 129000
$ws.Range("L71").Value = 131295
$ws.Range("M71").Value = -125256
$ws.Range("N71").Value = -138783
$ws.Range("H86").Value = 5879.524
$ws.Range("I86").Value = 5732.923
$ws.Range("K86").Value = 5732.923
$ws.Range("M86").Value = -4609.923
$ws.Range("H89").Value = 5879.524
$ws.Range("I89").Value = 5732.923
$ws.Range("K89").Value = 28664.615
$ws.Range("M89").Value = -23048.615
$ws.Range("H99").Value = 7038.778
$ws.Range("I99").Value = 6927.091
$ws.Range("K99").Value = 6927.091
$ws.Range("M99").Value = -5429.091
$ws.Range("H126").Value = 7038.778
$ws.Range("I126").Value = 6927.091
$ws.Range("K126").Value = 20781.273
$ws.Range("M126").Value = -18311.273
$ws.Range("H132").Value = 1178379.9
$ws.Range("I132").Value = 1430429.2
$ws.Range("K132").Value = 4291287.6
$ws.Range("M132").Value = -4288757.6
$ws.Range("H141").Value = 854452.4399999999
$ws.Range("J141").Value = 934897.7
$ws.Range("L141").Value = 934897.7
$ws.Range("N141").Value = -945257.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4110314.5
$ws.Range("J4").Value = 8488723
$ws.Range("L4").Value = 25466169
$ws.Range("N4").Value = -25466393
$ws.Range("H56").Value = 6990.125
$ws.Range("I56").Value = 6990.125
$ws.Range("K56").Value = 6990.125
$ws.Range("M56").Value = -6460.125
$ws.Range("H61").Value = 8375.5
$ws.Range("I61").Value = 84.666664
$ws.Range("J61").Value = 16666.334
$ws.Range("K61").Value = 253.999992
$ws.Range("L61").Value = 49999.00199999999
$ws.Range("M61").Value = -38.99999199999999
$ws.Range("N61").Value = -50429.00199999999
$ws.Range("H113").Value = 2402.0715
$ws.Range("J113").Value = 2694.25
$ws.Range("L113").Value = 8082.75
$ws.Range("N113").Value = -12422.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3529.4546
$ws.Range("I70").Value = 3455.75
$ws.Range("J70").Value = 3726
$ws.Range("K70").Value = 3455.75
$ws.Range("L70").Value = 3726
$ws.Range("M70").Value = -3185.75
$ws.Range("N70").Value = -4266
$ws.Range("H73").Value = 3529.4546
$ws.Range("I73").Value = 3455.75
$ws.Range("J73").Value = 3726
$ws.Range("K73").Value = 3455.75
$ws.Range("L73").Value = 3726
$ws.Range("M73").Value = -2519.75
$ws.Range("N73").Value = -5598
$ws.Range("H80").Value = 3579.8293
$ws.Range("I80").Value = 3319.3157
$ws.Range("J80").Value = 3804.818
$ws.Range("K80").Value = 3319.3157
$ws.Range("L80").Value = 3804.818
$ws.Range("M80").Value = -2321.3157
$ws.Range("N80").Value = -5800.818
$ws.Range("H83").Value = 3579.8293
$ws.Range("I83").Value = 3319.3157
$ws.Range("J83").Value = 3804.818
$ws.Range("K83").Value = 16596.5785
$ws.Range("L83").Value = 19024.09
$ws.Range("M83").Value = -11604.5785
$ws.Range("N83").Value = -29008.09

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8467.883
$ws.Range("J7").Value = 9622.625
$ws.Range("L7").Value = 9622.625
$ws.Range("N7").Value = -9846.625
$ws.Range("H40").Value = 6716.933
$ws.Range("I40").Value = 5972.3477
$ws.Range("J40").Value = 9163.429
$ws.Range("K40").Value = 5972.3477
$ws.Range("L40").Value = 9163.429
$ws.Range("M40").Value = -5836.3477
$ws.Range("N40").Value = -9435.429
$ws.Range("H61").Value = 65875.125
$ws.Range("I61").Value = 65875.125
$ws.Range("K61").Value = 65875.125
$ws.Range("M61").Value = -65673.125
$ws.Range("H74").Value = 24750
$ws.Range("I74").Value = 15000
$ws.Range("K74").Value = 15000
$ws.Range("M74").Value = -14002
$ws.Range("H77").Value = 24750
$ws.Range("I77").Value = 15000
$ws.Range("K77").Value = 45000
$ws.Range("M77").Value = -40008
$ws.Range("H113").Value = 65875.125
$ws.Range("I113").Value = 65875.125
$ws.Range("K113").Value = 65875.125
$ws.Range("M113").Value = -63705.125
$ws.Range("H122").Value = 6978.5
$ws.Range("I122").Value = 6575.7856
$ws.Range("J122").Value = 7683.25
$ws.Range("K122").Value = 19727.3568
$ws.Range("L122").Value = 23049.75
$ws.Range("M122").Value = -17277.3568
$ws.Range("N122").Value = -27949.75
$ws.Range("H126").Value = 8467.883
$ws.Range("J126").Value = 9622.625
$ws.Range("L126").Value = 28867.875
$ws.Range("N126").Value = -33807.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 30814.889
$ws.Range("I2").Value = 30814.889
$ws.Range("K2").Value = 30814.889
$ws.Range("M2").Value = -30702.889
$ws.Range("H126").Value = 2466.1292
$ws.Range("I126").Value = 1960.2174
$ws.Range("K126").Value = 5880.6522
$ws.Range("M126").Value = -3410.6522
$ws.Range("H132").Value = 1699.4595
$ws.Range("I132").Value = 1642.6364
$ws.Range("K132").Value = 4927.9092
$ws.Range("M132").Value = -2397.9092
$ws.Range("H136").Value = 4420.5483
$ws.Range("J136").Value = 2558.375
$ws.Range("L136").Value = 7675.125
$ws.Range("N136").Value = -12775.125
